$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.776.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.100.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.03'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.58'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.42%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.13%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.53'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.412.52'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.809'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.27%  '
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.51'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.103.99'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.816.17'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.06'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.76%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0840'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.68'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.66'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.53'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.98%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.66%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.34'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.47'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.34%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.59%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.54'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.41'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.57'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0228'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.54'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.533.78'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.76'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.24%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0910'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.14%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.96'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.295.51'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.28%  '
